# Auto-generated edit script: update leve-profit calc columns (H-N)
# across ALC/ARM/BSM/CUL/GSM/LTW/WVR sheets per scheduled-runner price refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 1049.375
$ws.Range("I38").Value = 58.75
$ws.Range("J38").Value = 2040
$ws.Range("K38").Value = 176.25
$ws.Range("L38").Value = 6120
$ws.Range("M38").Value = 195.75
$ws.Range("N38").Value = -6864
$ws.Range("H40").Value = 6960.222
$ws.Range("I40").Value = 4925.0713
$ws.Range("J40").Value = 14083.25
$ws.Range("K40").Value = 4925.0713
$ws.Range("L40").Value = 14083.25
$ws.Range("M40").Value = -4750.0713
$ws.Range("N40").Value = -14433.25
$ws.Range("H76").Value = 14712761
$ws.Range("J76").Value = 4295.2
$ws.Range("L76").Value = 4295.2
$ws.Range("N76").Value = -4925.2
$ws.Range("H79").Value = 14712761
$ws.Range("J79").Value = 4295.2
$ws.Range("L79").Value = 4295.2
$ws.Range("N79").Value = -6479.2
$ws.Range("H86").Value = 4518.2085
$ws.Range("I86").Value = 1694.9231
$ws.Range("J86").Value = 7854.8184
$ws.Range("K86").Value = 1694.9231
$ws.Range("L86").Value = 7854.8184
$ws.Range("M86").Value = -571.9231
$ws.Range("N86").Value = -10100.8184
$ws.Range("H89").Value = 4518.2085
$ws.Range("I89").Value = 1694.9231
$ws.Range("J89").Value = 7854.8184
$ws.Range("K89").Value = 8474.6155
$ws.Range("L89").Value = 39274.092
$ws.Range("M89").Value = -2858.6155
$ws.Range("N89").Value = -50506.092
$ws.Range("H92").Value = 3393.3125
$ws.Range("I92").Value = 8484.666999999999
$ws.Range("J92").Value = 338.5
$ws.Range("K92").Value = 8484.666999999999
$ws.Range("L92").Value = 338.5
$ws.Range("M92").Value = -7236.666999999999
$ws.Range("N92").Value = -2834.5
$ws.Range("H96").Value = 1146.8462
$ws.Range("I96").Value = 896.0952
$ws.Range("J96").Value = 2200
$ws.Range("K96").Value = 2688.2856
$ws.Range("L96").Value = 6600
$ws.Range("M96").Value = -1315.2856
$ws.Range("N96").Value = -9346
$ws.Range("H97").Value = 6000
$ws.Range("J97").Value = 6000
$ws.Range("L97").Value = 18000
$ws.Range("N97").Value = -18992

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H117").Value = 33372.5
$ws.Range("J117").Value = 33372.5
$ws.Range("L117").Value = 33372.5
$ws.Range("N117").Value = -42550.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 5420.8335
$ws.Range("I86").Value = 6581.75
$ws.Range("J86").Value = 3099
$ws.Range("K86").Value = 6581.75
$ws.Range("L86").Value = 3099
$ws.Range("M86").Value = -5458.75
$ws.Range("N86").Value = -5345
$ws.Range("H89").Value = 5420.8335
$ws.Range("I89").Value = 6581.75
$ws.Range("J89").Value = 3099
$ws.Range("K89").Value = 32908.75
$ws.Range("L89").Value = 15495
$ws.Range("M89").Value = -27292.75
$ws.Range("N89").Value = -26727

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H94").Value = 5826.125
$ws.Range("I94").Value = 3012
$ws.Range("J94").Value = 6228.143
$ws.Range("K94").Value = 9036
$ws.Range("L94").Value = 18684.429
$ws.Range("M94").Value = -8360
$ws.Range("N94").Value = -20036.429
$ws.Range("H95").Value = 7266.6665
$ws.Range("J95").Value = 7266.6665
$ws.Range("L95").Value = 21799.9995
$ws.Range("N95").Value = -25917.9995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5118.244
$ws.Range("I70").Value = 5573.643
$ws.Range("J70").Value = 4137.385
$ws.Range("K70").Value = 5573.643
$ws.Range("L70").Value = 4137.385
$ws.Range("M70").Value = -5303.643
$ws.Range("N70").Value = -4677.385
$ws.Range("H73").Value = 5118.244
$ws.Range("I73").Value = 5573.643
$ws.Range("J73").Value = 4137.385
$ws.Range("K73").Value = 5573.643
$ws.Range("L73").Value = 4137.385
$ws.Range("M73").Value = -4637.643
$ws.Range("N73").Value = -6009.385
$ws.Range("H80").Value = 4565.0415
$ws.Range("I80").Value = 4667.0454
$ws.Range("J80").Value = 3443
$ws.Range("K80").Value = 4667.0454
$ws.Range("L80").Value = 3443
$ws.Range("M80").Value = -3669.0454
$ws.Range("N80").Value = -5439
$ws.Range("H83").Value = 4565.0415
$ws.Range("I83").Value = 4667.0454
$ws.Range("J83").Value = 3443
$ws.Range("K83").Value = 23335.227
$ws.Range("L83").Value = 17215
$ws.Range("M83").Value = -18343.227
$ws.Range("N83").Value = -27199
$ws.Range("H113").Value = 38472452
$ws.Range("I113").Value = 62513616
$ws.Range("J113").Value = 6591.2
$ws.Range("K113").Value = 62513616
$ws.Range("L113").Value = 6591.2
$ws.Range("M113").Value = -62511446
$ws.Range("N113").Value = -10931.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2362.0334
$ws.Range("I68").Value = 2154.8333
$ws.Range("J68").Value = 2672.8333
$ws.Range("K68").Value = 2154.8333
$ws.Range("L68").Value = 2672.8333
$ws.Range("M68").Value = -1405.8333
$ws.Range("N68").Value = -4170.8333
$ws.Range("H71").Value = 2362.0334
$ws.Range("I71").Value = 2154.8333
$ws.Range("J71").Value = 2672.8333
$ws.Range("K71").Value = 10774.1665
$ws.Range("L71").Value = 13364.1665
$ws.Range("M71").Value = -7030.166499999999
$ws.Range("N71").Value = -20852.1665

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1845
$ws.Range("I81").Value = 866.36365
$ws.Range("J81").Value = 3382.8572
$ws.Range("K81").Value = 1732.7273
$ws.Range("L81").Value = 6765.7144
$ws.Range("M81").Value = -671.7273
$ws.Range("N81").Value = -8887.714400000001
$ws.Range("H84").Value = 1845
$ws.Range("I84").Value = 866.36365
$ws.Range("J84").Value = 3382.8572
$ws.Range("K84").Value = 8663.636500000001
$ws.Range("L84").Value = 33828.572
$ws.Range("M84").Value = -3359.636500000001
$ws.Range("N84").Value = -44436.572
$ws.Range("H96").Value = 1974.4615
$ws.Range("I96").Value = 1599.75
$ws.Range("J96").Value = 2574
$ws.Range("K96").Value = 1599.75
$ws.Range("L96").Value = 2574
$ws.Range("M96").Value = -226.75
$ws.Range("N96").Value = -5320
$ws.Range("H126").Value = 1367.68
$ws.Range("I126").Value = 1158.174
$ws.Range("J126").Value = 3777
$ws.Range("K126").Value = 3474.522
$ws.Range("L126").Value = 11331
$ws.Range("M126").Value = -1004.522
$ws.Range("N126").Value = -16271

